$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: rows 2 and 3
$wsZh.Range("E2").Value = "2016-03-17 07:01:47"
$wsZh.Range("H2").Value = "2016-03-17 07:02:28"
$wsZh.Range("E3").Value = "2016-03-17 07:01:47"
$wsZh.Range("H3").Value = "2016-03-17 07:02:28"

# de-de sheet: rows 2 and 3
$wsDe.Range("E2").Value = "2016-03-17 07:01:56"
$wsDe.Range("H2").Value = "2016-03-17 07:02:41"
$wsDe.Range("E3").Value = "2016-03-17 07:01:56"
$wsDe.Range("H3").Value = "2016-03-17 07:02:41"
